# Update the design doc (slide 3 / "DataCounting" -> "ProtocolCheker" area).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# --- Shape "TextBox 41" (id=42): reposition/resize + replace the 3 runs
#     ("D" / "ata" / " counting") with a single "ProtocolCheker" run. ---
$tb41 = $s.Shapes.Item(26)

# Move/resize (top & height stay the same).
$tb41.Left = 414.72023622047243
$tb41.Width = 144.8876090551181

# Replace the text; this collapses all runs into one run that inherits the
# first run's character formatting (bold + Consolas latin typeface).
$tb41.TextFrame.TextRange.Text = "ProtocolCheker"

# Align the run's language with the target ("en-SG") and clear any
# inherited effect list explicitly (writes an empty <a:effectLst/>).
$tb41.TextFrame.TextRange.LanguageID = "en-SG"
$tb41.TextFrame.TextRange.Font.Shadow = $false

# --- New connector: "Straight Arrow Connector 46" ---
# Duplicate an existing purple arrow connector so the new shape keeps the
# same line style/quick-style (p:style) and connector-lock metadata, then
# move it into place and flatten it to a horizontal segment.
$srcConn = $s.Shapes.Item(2)
$newConnRange = $srcConn.Duplicate()
$newConn = $newConnRange.Item(1)
$newConn.Name = "Straight Arrow Connector 46"
$newConn.Left = 559.6967716535434
$newConn.Top = 478.24697913385825
$newConn.Width = 29.12708661417323
$newConn.Height = 0
